# "se modificaron las cuentas de todos los datasources"
# Update the NroCuenta ("account number") value used by the Sheet1
# datasource row, clear the leftover one-off font formatting that used to
# highlight that cell, and leave the active selection parked on it - same
# as what the author's Excel session shows in the saved workbook.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")

# New account number for the Sheet1 datasource row.
$ws1.Range("E2").Value = 7635954411

# The old value carried a one-off "FF444444" font color style that singled
# this cell out; the refreshed value goes back to the sheet's normal,
# unstyled look.
$ws1.Range("E2").Style = "Normal"

# Leave the sheet with E2 (the cell that was just edited) selected/active,
# matching the saved cursor position.
$ws1.Activate()
$ws1.Range("E2").Select()

$wb.Save()
